$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9545.454
$ws.Range("J40").Value = 9545.454
$ws.Range("L40").Value = 9545.454
$ws.Range("N40").Value = -9895.454

$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10138

$ws.Range("H51").Value = 9999
$ws.Range("I51").Value = 9999
$ws.Range("K51").Value = 9999
$ws.Range("M51").Value = -9515

$ws.Range("H69").Value = 5333.3335
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16748

$ws.Range("H72").Value = 5333.3335
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53736

$ws.Range("H74").Value = 4750
$ws.Range("I74").Value = 4750
$ws.Range("K74").Value = 4750
$ws.Range("M74").Value = -3814

$ws.Range("H77").Value = 4750
$ws.Range("I77").Value = 4750
$ws.Range("K77").Value = 23750
$ws.Range("M77").Value = -19070

$ws.Range("H111").Value = 1954.1177
$ws.Range("I111").Value = 2168
$ws.Range("J111").Value = 1648.5714
$ws.Range("K111").Value = 6504
$ws.Range("L111").Value = 4945.7142
$ws.Range("M111").Value = -3437
$ws.Range("N111").Value = -11079.7142

$ws.Range("H116").Value = 7262.909
$ws.Range("J116").Value = 7585.5713
$ws.Range("L116").Value = 7585.5713
$ws.Range("N116").Value = -14469.5713

$ws.Range("H135").Value = 856.8570999999999
$ws.Range("I135").Value = 856.8570999999999
$ws.Range("K135").Value = 7711.7139
$ws.Range("M135").Value = -5176.7139

$ws.Range("H138").Value = 14366.611
$ws.Range("I138").Value = 7750
$ws.Range("J138").Value = 14755.823
$ws.Range("K138").Value = 23250
$ws.Range("L138").Value = 44267.469
$ws.Range("M138").Value = -18110
$ws.Range("N138").Value = -54547.469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2731.6
$ws.Range("I45").Value = 1602.909
$ws.Range("K45").Value = 1602.909
$ws.Range("M45").Value = -1225.909

$ws.Range("H132").Value = 2010.5454
$ws.Range("I132").Value = 1821.9333
$ws.Range("J132").Value = 2414.7144
$ws.Range("K132").Value = 5465.7999
$ws.Range("L132").Value = 7244.1432
$ws.Range("M132").Value = -2935.7999
$ws.Range("N132").Value = -12304.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 816.44446
$ws.Range("I80").Value = 500.25
$ws.Range("K80").Value = 500.25
$ws.Range("M80").Value = 497.75

$ws.Range("H83").Value = 816.44446
$ws.Range("I83").Value = 500.25
$ws.Range("K83").Value = 2501.25
$ws.Range("M83").Value = 2490.75

$ws.Range("H107").Value = 1140
$ws.Range("I107").Value = 1095.3334
$ws.Range("K107").Value = 1095.3334
$ws.Range("M107").Value = 824.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3602.4443
$ws.Range("I31").Value = 2506
$ws.Range("J31").Value = 4479.6
$ws.Range("K31").Value = 2506
$ws.Range("L31").Value = 4479.6
$ws.Range("M31").Value = -2211
$ws.Range("N31").Value = -5069.6

$ws.Range("H34").Value = 3602.4443
$ws.Range("I34").Value = 2506
$ws.Range("J34").Value = 4479.6
$ws.Range("K34").Value = 2506
$ws.Range("L34").Value = 4479.6
$ws.Range("M34").Value = -2304
$ws.Range("N34").Value = -4883.6

$ws.Range("H86").Value = 34853476
$ws.Range("J86").Value = 9982
$ws.Range("L86").Value = 9982
$ws.Range("N86").Value = -12228

$ws.Range("H89").Value = 34853476
$ws.Range("J89").Value = 9982
$ws.Range("L89").Value = 49910
$ws.Range("N89").Value = -61142

$ws.Range("H107").Value = 1449.6666
$ws.Range("I107").Value = 739.1111
$ws.Range("K107").Value = 739.1111
$ws.Range("M107").Value = 1180.8889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 518.5
$ws.Range("I8").Value = 518.5
$ws.Range("K8").Value = 1555.5
$ws.Range("M8").Value = -1416.5

$ws.Range("H12").Value = 442.75
$ws.Range("J12").Value = 455.63635
$ws.Range("L12").Value = 1366.90905
$ws.Range("N12").Value = -1712.90905

$ws.Range("H23").Value = 125
$ws.Range("I23").Value = 131.25
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 393.75
$ws.Range("L23").Value = 300
$ws.Range("M23").Value = -158.75
$ws.Range("N23").Value = -770

$ws.Range("H25").Value = 500.25
$ws.Range("I25").Value = 500.25
$ws.Range("K25").Value = 1500.75
$ws.Range("M25").Value = -1331.75

$ws.Range("H30").Value = 500.25
$ws.Range("I30").Value = 500.25
$ws.Range("K30").Value = 1500.75
$ws.Range("M30").Value = -1398.75

$ws.Range("H36").Value = 3027.8572
$ws.Range("I36").Value = 199.33333
$ws.Range("J36").Value = 19999
$ws.Range("K36").Value = 597.99999
$ws.Range("L36").Value = 59997
$ws.Range("M36").Value = -428.99999
$ws.Range("N36").Value = -60335

$ws.Range("H38").Value = 144.2
$ws.Range("I38").Value = 57.6
$ws.Range("J38").Value = 317.4
$ws.Range("K38").Value = 172.8
$ws.Range("L38").Value = 952.1999999999999
$ws.Range("M38").Value = 174.2
$ws.Range("N38").Value = -1646.2

$ws.Range("H80").Value = 1455.8
$ws.Range("J80").Value = 1039.5
$ws.Range("L80").Value = 3118.5
$ws.Range("N80").Value = -4990.5

$ws.Range("H83").Value = 1455.8
$ws.Range("J83").Value = 1039.5
$ws.Range("L83").Value = 9355.5
$ws.Range("N83").Value = -18715.5

$ws.Range("H86").Value = 7461.769
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 9400.299999999999
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 28200.9
$ws.Range("M86").Value = -1814
$ws.Range("N86").Value = -30572.9

$ws.Range("H89").Value = 7461.769
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 9400.299999999999
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 84602.7
$ws.Range("M89").Value = -3072
$ws.Range("N89").Value = -96458.7

$ws.Range("H94").Value = 750
$ws.Range("I94").Value = 750
$ws.Range("K94").Value = 2250
$ws.Range("M94").Value = -1574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7146564
$ws.Range("J70").Value = 4061.5
$ws.Range("L70").Value = 4061.5
$ws.Range("N70").Value = -4601.5

$ws.Range("H73").Value = 7146564
$ws.Range("J73").Value = 4061.5
$ws.Range("L73").Value = 4061.5
$ws.Range("N73").Value = -5933.5

$ws.Range("H92").Value = 19124.5
$ws.Range("J92").Value = 19124.5
$ws.Range("L92").Value = 19124.5
$ws.Range("N92").Value = -22868.5

$ws.Range("H113").Value = 1241.2858
$ws.Range("I113").Value = 1241.2858
$ws.Range("K113").Value = 1241.2858
$ws.Range("M113").Value = 928.7141999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9144.700000000001
$ws.Range("I132").Value = 6636
$ws.Range("J132").Value = 14998.333
$ws.Range("K132").Value = 19908
$ws.Range("L132").Value = 44994.999
$ws.Range("M132").Value = -17378
$ws.Range("N132").Value = -50054.999
